$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell B2: team -> comment label
$ws.Range("B2").Value = "コメント"

# Apply plain (non-highlighted) style to C4:G33 by copying format from an existing plain-style cell (H9)
$ws.Range("H9").Copy()
$ws.Range("C4:G33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update staff names (column A) and comments (column B)
$ws.Range("A4").Value = "スタッフ1"
$ws.Range("B4").ClearContents()

$ws.Range("A5").Value = "スタッフ2"
$ws.Range("B5").Value = "夜勤4回まで"

$ws.Range("A6").Value = "スタッフ3"
$ws.Range("B6").Value = "夜勤4回まで"

$ws.Range("A7").Value = "スタッフ4"
$ws.Range("B7").ClearContents()

$ws.Range("A8").Value = "スタッフ5"
$ws.Range("B8").ClearContents()

$ws.Range("A9").Value = "スタッフ6"
$ws.Range("B9").ClearContents()

$ws.Range("A10").Value = "スタッフ7"
$ws.Range("B10").ClearContents()

$ws.Range("A11").Value = "スタッフ8"
$ws.Range("B11").ClearContents()

$ws.Range("A12").Value = "スタッフ9"
$ws.Range("B12").Value = "夜勤3回まで"

$ws.Range("A13").Value = "スタッフ10"
$ws.Range("B13").ClearContents()

$ws.Range("A14").Value = "スタッフ11"
$ws.Range("B14").ClearContents()

$ws.Range("A15").Value = "スタッフ12"
$ws.Range("B15").ClearContents()

$ws.Range("A16").Value = "スタッフ13"
$ws.Range("B16").ClearContents()

$ws.Range("A17").Value = "スタッフ14"
$ws.Range("B17").ClearContents()

$ws.Range("A18").Value = "スタッフ15"
$ws.Range("B18").ClearContents()

$ws.Range("A19").Value = "スタッフ16"
$ws.Range("B19").Value = "新人　月前半長夜勤なし　"

$ws.Range("A20").Value = "スタッフ17"
$ws.Range("B20").ClearContents()

$ws.Range("A21").Value = "スタッフ18"
$ws.Range("B21").ClearContents()

$ws.Range("A22").Value = "スタッフ19"
$ws.Range("B22").ClearContents()

$ws.Range("A23").Value = "スタッフ20"
$ws.Range("B23").Value = "土日休み日勤のみ"

$ws.Range("A24").Value = "スタッフ21"
$ws.Range("B24").Value = "夜勤土日のみ3回まで"

$ws.Range("A25").Value = "スタッフ22"
$ws.Range("B25").Value = "長入明　水木金3回まで"

$ws.Range("A26").Value = "スタッフ23"
$ws.Range("B26").ClearContents()

$ws.Range("A27").Value = "スタッフ24"
$ws.Range("B27").ClearContents()

$ws.Range("A28").Value = "スタッフ25"
$ws.Range("B28").ClearContents()

$ws.Range("A29").Value = "スタッフ26"
$ws.Range("B29").ClearContents()

$ws.Range("A30").Value = "スタッフ27"
$ws.Range("B30").ClearContents()

$ws.Range("A31").Value = "スタッフ28"
$ws.Range("B31").ClearContents()

$ws.Range("A32").Value = "スタッフ29"
$ws.Range("B32").ClearContents()

$ws.Range("A33").Value = "スタッフ30"
$ws.Range("B33").Value = "新人　月前半長夜勤なし"
